# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh:
#   - re-ranks a few countries whose totals crossed each other
#     (Camerun/Honduras, Dominica/Fiyi, Santa Sede/Islas Turcas y Caicos,
#      Papua Nueva Guinea/Islas Virgenes Britanicas)
#   - refreshes case counts for a number of country rows
#   - refreshes the "datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow {
    param($Sheet, $Row, $Name, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)

    $Sheet.Range("A$Row").Value = $Name
    $Sheet.Range("B$Row").Value = $Total
    $Sheet.Range("C$Row").Value = $Nuevos
    $Sheet.Range("D$Row").Value = $Activos
    $Sheet.Range("E$Row").Value = $Recuperados
    $Sheet.Range("F$Row").Value = $Criticos
    $Sheet.Range("G$Row").Value = $MuertesHoy
    $Sheet.Range("H$Row").Value = $Muertes
}

# --- Rows with simple numeric refreshes (order / labels unchanged) ---
Set-CountryRow $ws 4   "Estados Unidos"       2330322 33132 972505 1235838 0 572 121979
Set-CountryRow $ws 5   "Brasil"               1070139 31571 543186 476895  0 968 50058
Set-CountryRow $ws 21  "Canada"               101019  390   63488  29121   0 64  8410
Set-CountryRow $ws 35  "Argentina"            41204   1634  12206  28006   0 13  992
Set-CountryRow $ws 47  "Panama"               25222   948   14359  10370   0 8   493
Set-CountryRow $ws 91  "Venezuela"            3789    198   835    2921    0 3   33
Set-CountryRow $ws 111 "Islandia"             1822    3     1805   7       0 0   10
Set-CountryRow $ws 128 "Republica de Chipre"  985     0     824    142     0 0   19
Set-CountryRow $ws 151 "Libia"                544     24    98     436     0 0   10
Set-CountryRow $ws 160 "Surinam"               303     10    74     221     0 0   8
Set-CountryRow $ws 173 "Eritrea"               143     1     39     104     0 0   0

# --- Re-ranked pair: Camerun overtakes Honduras ---
Set-CountryRow $ws 65 "Camerun"  11610 329 7702 3607 0 1 301
Set-CountryRow $ws 66 "Honduras" 11258 519 1214 9695 0 6 349

# --- Re-ranked pair: Fiyi overtakes Dominica (counts tied, only labels swap) ---
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# --- Re-ranked pair: Santa Sede overtakes Islas Turcas y Caicos ---
Set-CountryRow $ws 208 "Santa Sede"            12 0 12 0 0 0 0
Set-CountryRow $ws 209 "Islas Turcas y Caicos" 12 0 11 0 0 0 1

# --- Re-ranked pair: Papua Nueva Guinea overtakes Islas Virgenes Britanicas ---
Set-CountryRow $ws 213 "Papua Nueva Guinea"        8 0 8 0 0 0 0
Set-CountryRow $ws 214 "Islas Virgenes Britanicas" 8 0 7 0 0 0 1

# --- Timestamp refresh ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Junio de 2020 a las 02:23"
